$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 371.5357
$ws.Range("I53").Value = 399.1
$ws.Range("J53").Value = 302.625
$ws.Range("K53").Value = 399.1
$ws.Range("L53").Value = 302.625
$ws.Range("M53").Value = 237.9
$ws.Range("N53").Value = -1576.625
$ws.Range("H125").Value = 885.8333
$ws.Range("J125").Value = 1321.6666
$ws.Range("L125").Value = 11894.9994
$ws.Range("N125").Value = -16814.9994
$ws.Range("H129").Value = 3182.9285
$ws.Range("I129").Value = 3677.3333
$ws.Range("K129").Value = 11031.9999
$ws.Range("M129").Value = -6031.999899999999
$ws.Range("H132").Value = 61089.875
$ws.Range("I132").Value = 68064.42999999999
$ws.Range("J132").Value = 12268
$ws.Range("K132").Value = 204193.29
$ws.Range("L132").Value = 36804
$ws.Range("M132").Value = -201663.29
$ws.Range("N132").Value = -41864
$ws.Range("H135").Value = 408.46875
$ws.Range("I135").Value = 362.36667
$ws.Range("K135").Value = 3261.30003
$ws.Range("M135").Value = -726.3000299999999
$ws.Range("H137").Value = 2804.0454
$ws.Range("I137").Value = 2770.9
$ws.Range("J137").Value = 2831.6667
$ws.Range("K137").Value = 8312.700000000001
$ws.Range("L137").Value = 8495.000100000001
$ws.Range("M137").Value = -5762.700000000001
$ws.Range("N137").Value = -13595.0001
$ws.Range("H138").Value = 3778.9155
$ws.Range("I138").Value = 2974.2727
$ws.Range("J138").Value = 3926.4333
$ws.Range("K138").Value = 8922.8181
$ws.Range("L138").Value = 11779.2999
$ws.Range("M138").Value = -3782.8181
$ws.Range("N138").Value = -22059.2999
$ws.Range("H141").Value = 1152.4546
$ws.Range("I141").Value = 1178.7
$ws.Range("K141").Value = 3536.1
$ws.Range("M141").Value = 1643.9

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5496139.5
$ws.Range("I32").Value = 6757420.5
$ws.Range("K32").Value = 6757420.5
$ws.Range("M32").Value = -6757133.5
$ws.Range("H45").Value = 3108.25
$ws.Range("I45").Value = 3118.0908
$ws.Range("K45").Value = 3118.0908
$ws.Range("M45").Value = -2741.0908
$ws.Range("H60").Value = 25000
$ws.Range("I60").Value = 25000
$ws.Range("K60").Value = 25000
$ws.Range("M60").Value = -24267
$ws.Range("H74").Value = 2283.4102
$ws.Range("I74").Value = 1401.6875
$ws.Range("K74").Value = 1401.6875
$ws.Range("M74").Value = -527.6875
$ws.Range("H77").Value = 2283.4102
$ws.Range("I77").Value = 1401.6875
$ws.Range("K77").Value = 7008.4375
$ws.Range("M77").Value = -2640.4375
$ws.Range("H102").Value = 2993.7334
$ws.Range("I102").Value = 2993.7334
$ws.Range("K102").Value = 2993.7334
$ws.Range("M102").Value = -1371.7334
$ws.Range("H122").Value = 1372.76
$ws.Range("I122").Value = 1176.3
$ws.Range("J122").Value = 2158.6
$ws.Range("K122").Value = 3528.9
$ws.Range("L122").Value = 6475.799999999999
$ws.Range("M122").Value = -1078.9
$ws.Range("N122").Value = -11375.8
$ws.Range("H132").Value = 1730.942
$ws.Range("I132").Value = 1522.6984
$ws.Range("K132").Value = 4568.0952
$ws.Range("M132").Value = -2038.0952

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1733.3334
$ws.Range("J7").Value = 1733.3334
$ws.Range("L7").Value = 1733.3334
$ws.Range("N7").Value = -1959.3334
$ws.Range("H86").Value = 2042
$ws.Range("I86").Value = 1681.3
$ws.Range("J86").Value = 2557.2856
$ws.Range("K86").Value = 1681.3
$ws.Range("L86").Value = 2557.2856
$ws.Range("M86").Value = -558.3
$ws.Range("N86").Value = -4803.2856
$ws.Range("H89").Value = 2042
$ws.Range("I89").Value = 1681.3
$ws.Range("J89").Value = 2557.2856
$ws.Range("K89").Value = 8406.5
$ws.Range("L89").Value = 12786.428
$ws.Range("M89").Value = -2790.5
$ws.Range("N89").Value = -24018.428
$ws.Range("H132").Value = 100495
$ws.Range("J132").Value = 100495
$ws.Range("L132").Value = 100495
$ws.Range("N132").Value = -110615

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2398
$ws.Range("I3").Value = 2764
$ws.Range("J3").Value = 1300
$ws.Range("K3").Value = 2764
$ws.Range("L3").Value = 1300
$ws.Range("M3").Value = -2651
$ws.Range("N3").Value = -1526
$ws.Range("H16").Value = 1642.6
$ws.Range("J16").Value = 2213
$ws.Range("L16").Value = 2213
$ws.Range("N16").Value = -2787
$ws.Range("H63").Value = 95203.25
$ws.Range("J63").Value = 100271
$ws.Range("L63").Value = 100271
$ws.Range("N63").Value = -101643
$ws.Range("H66").Value = 95203.25
$ws.Range("J66").Value = 100271
$ws.Range("L66").Value = 300813
$ws.Range("N66").Value = -307677
$ws.Range("H107").Value = 514.1852
$ws.Range("I107").Value = 392.5625
$ws.Range("K107").Value = 392.5625
$ws.Range("M107").Value = 1527.4375
$ws.Range("H113").Value = 1642.6
$ws.Range("J113").Value = 2213
$ws.Range("L113").Value = 2213
$ws.Range("N113").Value = -6553
$ws.Range("H122").Value = 3617.1875
$ws.Range("I122").Value = 2487.9
$ws.Range("K122").Value = 7463.700000000001
$ws.Range("M122").Value = -5013.700000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 262.41666
$ws.Range("J98").Value = 260.63635
$ws.Range("L98").Value = 781.90905
$ws.Range("N98").Value = -3777.90905
$ws.Range("H131").Value = 1713.931
$ws.Range("I131").Value = 849.5
$ws.Range("J131").Value = 1777.963
$ws.Range("K131").Value = 2548.5
$ws.Range("L131").Value = 5333.889
$ws.Range("M131").Value = 2491.5
$ws.Range("N131").Value = -15413.889
$ws.Range("H132").Value = 3825.6667
$ws.Range("I132").Value = 3100
$ws.Range("J132").Value = 4089.5454
$ws.Range("K132").Value = 27900
$ws.Range("L132").Value = 36805.9086
$ws.Range("M132").Value = -25370
$ws.Range("N132").Value = -41865.9086
$ws.Range("H140").Value = 2779.6086
$ws.Range("I140").Value = 2131.0667
$ws.Range("K140").Value = 6393.2001
$ws.Range("M140").Value = -1213.2001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 22500.75
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 22500.75
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").Value = 22500.75
$ws.Range("N40").Value = -22802.75
$ws.Range("H43").Value = 1935
$ws.Range("I43").Value = 1935
$ws.Range("K43").Value = 1935
$ws.Range("M43").Value = -1784
$ws.Range("H70").Value = 6039.1113
$ws.Range("I70").Value = 5955.905
$ws.Range("J70").Value = 6330.3335
$ws.Range("K70").Value = 5955.905
$ws.Range("L70").Value = 6330.3335
$ws.Range("M70").Value = -5685.905
$ws.Range("N70").Value = -6870.3335
$ws.Range("H73").Value = 6039.1113
$ws.Range("I73").Value = 5955.905
$ws.Range("J73").Value = 6330.3335
$ws.Range("K73").Value = 5955.905
$ws.Range("L73").Value = 6330.3335
$ws.Range("M73").Value = -5019.905
$ws.Range("N73").Value = -8202.333500000001
$ws.Range("H122").Value = 44918.36
$ws.Range("I122").Value = 61211.41
$ws.Range("J122").Value = 10295.625
$ws.Range("K122").Value = 183634.23
$ws.Range("L122").Value = 30886.875
$ws.Range("M122").Value = -181184.23
$ws.Range("N122").Value = -35786.875
$ws.Range("H126").Value = 3433
$ws.Range("I126").Value = 2878.3635
$ws.Range("K126").Value = 8635.0905
$ws.Range("M126").Value = -6165.0905
$ws.Range("H132").Value = 11824.976
$ws.Range("I132").Value = 10604.3
$ws.Range("J132").Value = 15154.091
$ws.Range("K132").Value = 31812.9
$ws.Range("L132").Value = 45462.273
$ws.Range("M132").Value = -29282.9
$ws.Range("N132").Value = -50522.273

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1249
$ws.Range("I22").Value = 748
$ws.Range("J22").Value = 1499.5
$ws.Range("K22").Value = 748
$ws.Range("L22").Value = 1499.5
$ws.Range("M22").Value = -453
$ws.Range("N22").Value = -2089.5
$ws.Range("H27").Value = 1249
$ws.Range("I27").Value = 748
$ws.Range("J27").Value = 1499.5
$ws.Range("K27").Value = 748
$ws.Range("L27").Value = 1499.5
$ws.Range("M27").Value = -641
$ws.Range("N27").Value = -1713.5
$ws.Range("H82").Value = 1710
$ws.Range("I82").Value = 1376
$ws.Range("K82").Value = 1376
$ws.Range("M82").Value = -1015
$ws.Range("H85").Value = 1710
$ws.Range("I85").Value = 1376
$ws.Range("K85").Value = 1376
$ws.Range("M85").Value = -128
$ws.Range("H93").Value = 3466.5
$ws.Range("I93").Value = 3466.5
$ws.Range("K93").Value = 3466.5
$ws.Range("M93").Value = -2218.5
$ws.Range("H122").Value = 4134.8
$ws.Range("I122").Value = 2918.625
$ws.Range("K122").Value = 8755.875
$ws.Range("M122").Value = -6305.875
$ws.Range("H132").Value = 2343.3125
$ws.Range("I132").Value = 2391.8076
$ws.Range("K132").Value = 7175.4228
$ws.Range("M132").Value = -4645.4228

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1867.4348
$ws.Range("I122").Value = 1647.3889
$ws.Range("J122").Value = 2659.6
$ws.Range("K122").Value = 4942.1667
$ws.Range("L122").Value = 7978.799999999999
$ws.Range("M122").Value = -2492.1667
$ws.Range("N122").Value = -12878.8
$ws.Range("H132").Value = 1802.5428
$ws.Range("I132").Value = 1534.0625
$ws.Range("K132").Value = 4602.1875
$ws.Range("M132").Value = -2072.1875
$ws.Range("H136").Value = 13456.694
$ws.Range("I136").Value = 3699.9092
$ws.Range("K136").Value = 11099.7276
$ws.Range("M136").Value = -8549.7276

Write-Output "Applied all Excalibur_Profits market-data refresh updates."